# Apply the edits described by the diff to the expo location-details sheet.
# The sheet's `name`/`rating`/`user_ratings_total` rows are kept alphabetically
# sorted by venue name. 'Grand Park Sports Campus' was removed and 'Indy Displays'
# was inserted in its alphabetical slot, so every row between them shifts up by one.
# A few unrelated review-count bumps are also applied (E11, E21, E35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = 473
$ws.Range("A17").Value = 16
$ws.Range("C17").Value = "Hamilton County Fairgrounds"
$ws.Range("D17").Value = 4.6
$ws.Range("E17").Value = 68
$ws.Range("A18").Value = 28
$ws.Range("C18").Value = "Hancock County Fairgrounds - Indiana"
$ws.Range("D18").Value = 4.2
$ws.Range("E18").Value = 333
$ws.Range("A19").Value = 18
$ws.Range("C19").Value = "Harvest Pavillion"
$ws.Range("D19").Value = 4.6
$ws.Range("E19").Value = 18
$ws.Range("A20").Value = 2
$ws.Range("C20").Value = "Indiana Black Expo Inc"
$ws.Range("D20").Value = 4.3
$ws.Range("E20").Value = 39
$ws.Range("A21").Value = 21
$ws.Range("C21").Value = "Indiana Convention Center"
$ws.Range("D21").Value = 4.5
$ws.Range("E21").Value = 528
$ws.Range("A22").Value = 20
$ws.Range("C22").Value = "Indiana Flower & Patio Show"
$ws.Range("D22").Value = 4.3
$ws.Range("E22").Value = 60
$ws.Range("A23").Value = 4
$ws.Range("C23").Value = "Indiana Latino Expo"
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("A24").Value = 15
$ws.Range("C24").Value = "Indiana State Fairgrounds & Event Center"
$ws.Range("D24").Value = 4.4
$ws.Range("E24").Value = 1344
$ws.Range("A25").Value = 32
$ws.Range("C25").Value = "Indiana State Numismatic Association"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("A26").Value = 26
$ws.Range("C26").Value = "Indianapolis Auto Show"
$ws.Range("D26").Value = 3.1
$ws.Range("E26").Value = 51
$ws.Range("A27").Value = 10
$ws.Range("C27").Value = "Indianapolis Chapter of Indiana Black Expo, Inc."
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("A28").Value = 17
$ws.Range("C28").Value = "Indianapolis Motor Speedway"
$ws.Range("D28").Value = 4.8
$ws.Range("E28").Value = 11013
$ws.Range("A29").Value = 9
$ws.Range("C29").Value = "Indy Air Expo"
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("A30").Value = 29
$ws.Range("C30").Value = "Indy Displays"
$ws.Range("D30").Value = 4.5
$ws.Range("E30").Value = 30
$ws.Range("A31").Value = 19
$ws.Range("A35").Value = 33
$ws.Range("E35").Value = 1757
$ws.Range("A37").Value = 30
$ws.Range("A40").Value = 31
